$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 220, shifting existing rows 220-315 down to 221-316
$ws.Rows("220:220").Insert()

# Populate the newly inserted row 220 with data
$ws.Cells.Item(220, 1).Value = 5
$ws.Cells.Item(220, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(220, 3).Value = "Maule"
$ws.Cells.Item(220, 4).Value = 44825
$ws.Cells.Item(220, 5).Value = 7
$ws.Cells.Item(220, 6).Value = 100112009
$ws.Cells.Item(220, 7).Value = "Acelga"
$ws.Cells.Item(220, 8).Value = "Sin especificar"
$ws.Cells.Item(220, 9).Value = "Primera"
$ws.Cells.Item(220, 10).Value = 400
$ws.Cells.Item(220, 11).Value = 2500
$ws.Cells.Item(220, 12).Value = 2500
$ws.Cells.Item(220, 13).Value = 2500
$ws.Cells.Item(220, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(220, 15).Value = "Región del Maule"
$ws.Cells.Item(220, 16).Value = 625
$ws.Cells.Item(220, 17).Value = 4
$ws.Cells.Item(220, 18).Value = "Hortaliza"

Write-Output "Row inserted and populated; dimension should now be A1:R316"
